$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.927.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.61%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.812.00"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.53%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.06"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4983"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.00%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3914"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.84%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09825"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +25.45%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.099"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.12%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.96"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.35%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.394"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.01%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.46"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.51%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.002"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.01%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.810.08"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.54%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.258"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.82%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001140"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +6.03%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.26"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06650"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.98%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.04%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.15"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.97%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.918"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.31%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.994.01"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.60%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.255"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.46%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.53"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.26%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.022.35"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.94%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.54"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.82%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.388"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.07%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.52"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.63%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1063"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.15%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.030"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.05%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.549"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.37%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.599"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.97%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06714"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.24%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.893"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.61%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.15%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2137"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.86%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.924"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.00%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.24"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.92%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6171"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.57%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.170"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.11%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.01%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.18"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.75%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5887"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.11%  "

# Row 46
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.692"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.03%  "

# Row 47
$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.281"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.94%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.47"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.03%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.930"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.64%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.177"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.70%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06773"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.72%  "
